# Mise a jour de l'initdata en vue de l'integration des donnees
#
# - Remove the 3 "refuge" rabbits (rows 8-10: Nutella, Defy, Fran)
# - Replace the numeric "proprietaire" ids (col M) with real usernames
# - Give every remaining animal a " " placeholder description where missing (col N)
# - Add a "date" column (O) with a per-row date value
# - Add "montant" / "date_sterilisation" / "poids" / "sante" columns (P:S), only
#   "montant" (P) carries data for now
# - Move the active selection to S2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- drop the three REFUGE/ABANDON/FOURRIERE rows --------------------------
$ws.Rows("8:10").Delete()

# --- proprietaire (M) becomes a username instead of a numeric id -----------
# the id-looking numbers carried the old "number" cell style (border etc.);
# the new usernames are plain, unstyled text, so clear formatting first.
$ws.Range("M2:M7").ClearFormats()
$ws.Range("M2").Value = "perreautbis.clementine"
$ws.Range("M3").Value = "callo.laurie"
$ws.Range("M4").Value = "chaplin.charles"
$ws.Range("M5").Value = "machado.chloe"
$ws.Range("M6").Value = "gauger.sabine"
$ws.Range("M7").Value = "serie.dexter"

# --- description (N) : rows that had no description get a single space ----
$ws.Range("N4").Value = " "
$ws.Range("N5").Value = " "
$ws.Range("N6").Value = " "
$ws.Range("N7").Value = " "

# --- date (O) : header keeps its old look, now stored as text (numFmt 49) --
$ws.Range("O1").Value = "date"
$ws.Range("O1").NumberFormat = "@"

$ws.Range("O2:O7").NumberFormat = "@"
$ws.Range("O2").Value = "2019-06-10"
$ws.Range("O3").Value = "2019-05-10"
$ws.Range("O4").Value = "2020-02-10"
$ws.Range("O5").Value = "2020-03-10"
$ws.Range("O6").Value = "2020-02-17"
$ws.Range("O7").Value = "2020-02-10"

# --- new trailing columns: montant / date_sterilisation / poids / sante ----
# Style them like the other header cells (font/fill/vertical-centre) but
# without the heavy table border, since they sit outside the original table.
$ws.Range("N1").Copy()
$ws.Range("P1:S1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1:S1").Borders.LineStyle = -4142  # xlLineStyleNone
$excel.CutCopyMode = $false

$ws.Range("P1").Value = "montant"
$ws.Range("Q1").Value = "date_sterilisation"
$ws.Range("R1").Value = "poids"
$ws.Range("S1").Value = "sante"

$ws.Range("P2").Value = 100
$ws.Range("P3").Value = 100
$ws.Range("P4").Value = 100
$ws.Range("P5").Value = 70
$ws.Range("P6").Value = 10
$ws.Range("P7").Value = 80

# --- restore a sane viewport / selection ------------------------------------
$ws.Range("S2").Select()
